# Edit script: bump stack-trace line numbers / stack frames in the
# asStyleTextNotExistingStyle expected-generation fixture (3.1.1 -> 3.2.0).
$d = $word.ActiveDocument

# M2DocEvaluator.caseQuery(M2DocEvaluator.java:587) -> M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)
$oldText = "M2DocEvaluator.caseQuery(M2DocEvaluator.java:587)"
$newText = "M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: M2DocEvaluator.caseQuery(M2DocEvaluator.java:587) -> M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)" }

# doSwitch(M2DocEvaluator.java:1242) -> doSwitch(M2DocEvaluator.java:1331)
$oldText = "doSwitch(M2DocEvaluator.java:1242)"
$newText = "doSwitch(M2DocEvaluator.java:1331)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: doSwitch(M2DocEvaluator.java:1242) -> doSwitch(M2DocEvaluator.java:1331)" }

# caseBlock(M2DocEvaluator.java:1467) -> caseBlock(M2DocEvaluator.java:1556)
$oldText = "caseBlock(M2DocEvaluator.java:1467)"
$newText = "caseBlock(M2DocEvaluator.java:1556)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: caseBlock(M2DocEvaluator.java:1467) -> caseBlock(M2DocEvaluator.java:1556)" }

# caseDocumentTemplate(M2DocEvaluator.java:297) -> caseDocumentTemplate(M2DocEvaluator.java:301)
$oldText = "caseDocumentTemplate(M2DocEvaluator.java:297)"
$newText = "caseDocumentTemplate(M2DocEvaluator.java:301)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: caseDocumentTemplate(M2DocEvaluator.java:297) -> caseDocumentTemplate(M2DocEvaluator.java:301)" }

# M2DocEvaluator.generate(M2DocEvaluator.java:282) -> M2DocEvaluator.generate(M2DocEvaluator.java:286)
$oldText = "M2DocEvaluator.generate(M2DocEvaluator.java:282)"
$newText = "M2DocEvaluator.generate(M2DocEvaluator.java:286)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: M2DocEvaluator.generate(M2DocEvaluator.java:282) -> M2DocEvaluator.generate(M2DocEvaluator.java:286)" }

# M2DocUtils.generate(M2DocUtils.java:845) -> M2DocUtils.generate(M2DocUtils.java:853)
$oldText = "M2DocUtils.generate(M2DocUtils.java:845)"
$newText = "M2DocUtils.generate(M2DocUtils.java:853)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: M2DocUtils.generate(M2DocUtils.java:845) -> M2DocUtils.generate(M2DocUtils.java:853)" }

# prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514) -> prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)
$oldText = "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)"
$newText = "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514) -> prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)" }

# AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421) -> AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)
$oldText = "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)"
$newText = "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421) -> AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)" }

# GeneratedMethodAccessor73.invoke(Unknown Source) -> GeneratedMethodAccessor5.invoke(Unknown Source)
$oldText = "GeneratedMethodAccessor73.invoke(Unknown Source)"
$newText = "GeneratedMethodAccessor5.invoke(Unknown Source)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: GeneratedMethodAccessor73.invoke(Unknown Source) -> GeneratedMethodAccessor5.invoke(Unknown Source)" }

# replace JDT/Eclipse test-runner trailer with Maven/Tycho/Equinox trailer
$oldText = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)" + [char]10 + "	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)" + [char]10 + "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)" + [char]10 + "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)" + [char]10 + "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)" + [char]10 + "	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$newText = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)" + [char]10 + "	at org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)" + [char]10 + "	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)" + [char]10 + "	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + "	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + "	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)" + [char]10 + "	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)" + [char]10 + "	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)" + [char]10 + "	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)" + [char]10 + "	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + "	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + "	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)" + [char]10 + "	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)" + [char]10 + "	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)" + [char]10 + "	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)" + [char]10 + "	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)" + [char]10 + "	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + "	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + "	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + "	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + "	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)" + [char]10 + "	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)" + [char]10 + "	at org.eclipse.equinox.launcher.Main.run(Main.java:1447)" + [char]10 + "	at org.eclipse.equinox.launcher.Main.main(Main.java:1420)"
$found = $d.Content.Find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) { Write-Host "NOT FOUND: replace JDT/Eclipse test-runner trailer with Maven/Tycho/Equinox trailer" }

Write-Host "Done"
